$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("rallies")

# Update existing last row (row 69): rally_no changes from 6 to 7
$ws.Cells.Item(69, 4).Value = 7

# Append a new row (row 70) describing a new rally
$ws.Cells.Item(70, 1).Value = 69      # A70 rally_id
$ws.Cells.Item(70, 2).Value = 1       # B70 match_id
$ws.Cells.Item(70, 3).Value = 3       # C70 set_number
$ws.Cells.Item(70, 4).Value = 7       # D70 rally_no
$ws.Cells.Item(70, 5).Value = "NOS"   # E70 side
$ws.Cells.Item(70, 6).Value = ""      # F70 position
$ws.Cells.Item(70, 7).Value = 4       # G70 player_number
$ws.Cells.Item(70, 8).Value = "MEIO"  # H70 action
$ws.Cells.Item(70, 9).Value = "PONTO" # I70 result
$ws.Cells.Item(70, 10).Value = "NOS"  # J70 who_scored
$ws.Cells.Item(70, 11).Value = 7      # K70 score_home
$ws.Cells.Item(70, 12).Value = 0      # L70 score_away
$ws.Cells.Item(70, 13).Value = "1 4 m"   # M70 raw_text
$ws.Cells.Item(70, 14).Value = "FRENTE"  # N70 position_zone
$ws.Cells.Item(70, 15).Value = "FRENTE"  # O70 pos_fb
$ws.Cells.Item(70, 16).Value = "FRENTE"  # P70 frente_fundo
